$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Drop the hidden `_xlchart.v1.*` defined names (chart-tracking helper
#    names Excel minted for the old layout; the rebuilt chart no longer
#    needs them).
# ---------------------------------------------------------------------------
$oldNames = @(
    "_xlchart.v1.0",
    "_xlchart.v1.1",
    "_xlchart.v1.2",
    "_xlchart.v1.3",
    "_xlchart.v1.4",
    "_xlchart.v1.5",
    "_xlchart.v1.6"
)
foreach ($nm in $oldNames) {
    $wb.Names.Item($nm).Delete()
}

# ---------------------------------------------------------------------------
# 2) Rebuild the data grid. The whole table shifts one column to the right
#    (A:J -> B:K) and gains a results-history layout: Bounds / Model /
#    Optimized 1 / Error / Optimized 2 / Error / Optimized 3 / Error.
#    Simplest & safest: wipe the old A1:J4 block, then write every cell of
#    the new layout explicitly.
# ---------------------------------------------------------------------------
$ws.Range("A1:J4").ClearContents() | Out-Null

# Header row (row 1, B1:K1) -- muscle names
$ws.Range("B1").Value = "bflh_r"
$ws.Range("C1").Value = "bfsh_r"
$ws.Range("D1").Value = "gaslat_r"
$ws.Range("E1").Value = "gasmed_r"
$ws.Range("F1").Value = "sart_r"
$ws.Range("G1").Value = "semimem_r"
$ws.Range("H1").Value = "semiten_r"
$ws.Range("I1").Value = "vasint_r"
$ws.Range("J1").Value = "vaslat_r"
$ws.Range("K1").Value = "vasmed_r"

# Row 2 -- Bounds
$ws.Range("A2").Value = "Bounds"
$ws.Range("E2").Value = 0.494442
$ws.Range("F2").Value = 0.436769
$ws.Range("G2").Value = 0.416913
$ws.Range("H2").Value = 0.471108
$ws.Range("I2").Value = 0.279678
$ws.Range("J2").Value = 0.292679
$ws.Range("K2").Value = 0.268264

# Row 3 -- Model
$ws.Range("A3").Value = "Model"
$ws.Range("B3").Value = 0.3179
$ws.Range("C3").Value = 0.104
$ws.Range("D3").Value = 0.432
$ws.Range("E3").Value = 0.457
$ws.Range("F3").Value = 0.124
$ws.Range("G3").Value = 0.33
$ws.Range("H3").Value = 0.245
$ws.Range("I3").Value = 0.2
$ws.Range("J3").Value = 0.2
$ws.Range("K3").Value = 0.198

# Row 4 -- Optimized 1
$ws.Range("A4").Value = "Optimized 1"
$ws.Range("B4").Value = 0.308030253144893
$ws.Range("C4").Value = 0.0472194466338797
$ws.Range("D4").Value = 0.416094636029345
$ws.Range("E4").Value = 0.43916405262659
$ws.Range("F4").Value = 0.482455634172036
$ws.Range("G4").Value = 0.323416426183307
$ws.Range("H4").Value = 0.229396111076588
$ws.Range("I4").Value = 0.319631097673541
$ws.Range("J4").Value = 0.190862802110195
$ws.Range("K4").Value = 0.181274563686359

# Row 5 -- Error vs Model: (Optimized1 - Model) / Model * 100
$ws.Range("A5").Value = "Error"
$ws.Range("B5").Formula = "=(B4-B3)/B3*100"
$ws.Range("C5").Formula = "=(C4-C3)/C3*100"
$ws.Range("D5").Formula = "=(D4-D3)/D3*100"
$ws.Range("E5").Formula = "=(E4-E3)/E3*100"
$ws.Range("F5").Formula = "=(F4-F3)/F3*100"
$ws.Range("G5").Formula = "=(G4-G3)/G3*100"
$ws.Range("H5").Formula = "=(H4-H3)/H3*100"
$ws.Range("I5").Formula = "=(I4-I3)/I3*100"
$ws.Range("J5").Formula = "=(J4-J3)/J3*100"
$ws.Range("K5").Formula = "=(K4-K3)/K3*100"

# Row 6 -- Optimized 2
$ws.Range("A6").Value = "Optimized 2"
$ws.Range("G6").Value = 0.338143335138717
$ws.Range("H6").Value = 0.245487214088891
$ws.Range("I6").Value = 0.198664661813162
$ws.Range("J6").Value = 0.2203206165746
$ws.Range("K6").Value = 0.196942178632675
$ws.Range("M6").Value = 0.198664661813162
$ws.Range("N6").Value = 0.2203206165746
$ws.Range("O6").Value = 0.196942178632675

# Row 7 -- Error vs Optimized 2 (B:F divide by blank -> #DIV/0!, G:K use row 4)
$ws.Range("A7").Value = "Error"
$ws.Range("B7").Formula = "=(B5-B6)/B6*100"
$ws.Range("C7").Formula = "=(C5-C6)/C6*100"
$ws.Range("D7").Formula = "=(D5-D6)/D6*100"
$ws.Range("E7").Formula = "=(E5-E6)/E6*100"
$ws.Range("F7").Formula = "=(F5-F6)/F6*100"
$ws.Range("G7").Formula = "=(G4-G6)/G6*100"
$ws.Range("H7").Formula = "=(H4-H6)/H6*100"
$ws.Range("I7").Formula = "=(I4-I6)/I6*100"
$ws.Range("J7").Formula = "=(J4-J6)/J6*100"
$ws.Range("K7").Formula = "=(K4-K6)/K6*100"

# Row 8 -- Optimized 3
$ws.Range("A8").Value = "Optimized 3"
$ws.Range("E8").Value = 0.456889766862435
$ws.Range("F8").Value = 0.279706945394611
$ws.Range("G8").Value = 0.338151646897526
$ws.Range("H8").Value = 0.245558047314078
$ws.Range("I8").Value = 0.196191019912824
$ws.Range("J8").Value = 0.211053926951964
$ws.Range("K8").Value = 0.1932993647691

# Row 9 -- Error vs Model (absolute $3 row reference)
$ws.Range("A9").Value = "Error"
$ws.Range("B9").Formula = "=(B8-B`$3)/B`$3*100"
$ws.Range("C9").Formula = "=(C8-C`$3)/C`$3*100"
$ws.Range("D9").Formula = "=(D8-D`$3)/D`$3*100"
$ws.Range("E9").Formula = "=(E8-E`$3)/E`$3*100"
$ws.Range("F9").Formula = "=(F8-F`$3)/F`$3*100"
$ws.Range("G9").Formula = "=(G8-G`$3)/G`$3*100"
$ws.Range("H9").Formula = "=(H8-H`$3)/H`$3*100"
$ws.Range("I9").Formula = "=(I8-I`$3)/I`$3*100"
$ws.Range("J9").Formula = "=(J8-J`$3)/J`$3*100"
$ws.Range("K9").Formula = "=(K8-K`$3)/K`$3*100"

# ---------------------------------------------------------------------------
# 3) Point the bar chart's category/value series at the new table location
#    (Sheet1!$A$1:$J$1 / Sheet1!$A$4:$J$4 -> $B$1:$K$1 / $B$5:$K$5) and shift
#    the chart itself one column right / one row taller to match.
#    (Do this before the column-A AutoFit below -- reordering it after the
#    width change perturbs the engine's px/EMU conversion for the anchor.)
# ---------------------------------------------------------------------------
$chartObj = $ws.ChartObjects().Item(1)
$series = $chartObj.Chart.SeriesCollection(1)
$series.Formula = "=SERIES(,Sheet1!`$B`$1:`$K`$1,Sheet1!`$B`$5:`$K`$5,1)"

$chartObj.Left = 686.3124212598425
$chartObj.Top = 9.37496062992126
$chartObj.Width = 700.5
$chartObj.Height = 231.0

# ---------------------------------------------------------------------------
# 4) Column A now holds the row labels ("Bounds"/"Model"/...) -- autofit it.
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# 5) Restore the active-cell selection to match where the editor left off.
# ---------------------------------------------------------------------------
$ws.Range("K16").Select() | Out-Null
